$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.48%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.49%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.691"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.70%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08161"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.74%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.054"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'7.13%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.759"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.23%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.541"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.75%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.949"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.24%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9201"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.31%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'1.25%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1942"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.77%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09435"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.82%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03693"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'5.34%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'10.22%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'-0.04%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006168"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.37%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'2.50%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-2.09%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.297"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-5.01%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1394"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.74%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'10.17%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04443"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.41%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001269"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.67%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.56%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001183"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'3.76%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02757"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'14.31%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05450"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.28%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007665"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.08%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009465"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.94%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1417"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.79%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002125"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.25%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01216"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'8.42%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006882"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.05%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'60.58%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003538"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'17.69%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.21%"
$ws.Range("E51").Style = "Normal"
